# Weekly update: insert two new price records (rows 291-292) for
# "Terminal Hortofrutícola Agro Chillán" / Repollo, shifting the existing
# rows 291:333 down to 293:335.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 291-292 (existing rows 291+ shift down by 2)
$ws.Rows("291:292").Insert()

# New row 291
$ws.Cells.Item(291, 1).Value = 7
$ws.Cells.Item(291, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(291, 3).Value = "Ñuble"
$ws.Cells.Item(291, 4).Value = 44984
$ws.Cells.Item(291, 5).Value = 16
$ws.Cells.Item(291, 6).Value = 100112006
$ws.Cells.Item(291, 7).Value = "Repollo"
$ws.Cells.Item(291, 8).Value = "Crespo record"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 500
$ws.Cells.Item(291, 11).Value = 1100
$ws.Cells.Item(291, 12).Value = 1200
$ws.Cells.Item(291, 13).Value = 1150
$ws.Cells.Item(291, 14).Value = "`$/unidad"
$ws.Cells.Item(291, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(291, 16).Value = 1150
$ws.Cells.Item(291, 17).Value = 1
$ws.Cells.Item(291, 18).Value = "Hortaliza"

# New row 292
$ws.Cells.Item(292, 1).Value = 7
$ws.Cells.Item(292, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(292, 3).Value = "Ñuble"
$ws.Cells.Item(292, 4).Value = 44984
$ws.Cells.Item(292, 5).Value = 16
$ws.Cells.Item(292, 6).Value = 100112006
$ws.Cells.Item(292, 7).Value = "Repollo"
$ws.Cells.Item(292, 8).Value = "Crespo record"
$ws.Cells.Item(292, 9).Value = "Segunda"
$ws.Cells.Item(292, 10).Value = 400
$ws.Cells.Item(292, 11).Value = 900
$ws.Cells.Item(292, 12).Value = 900
$ws.Cells.Item(292, 13).Value = 900
$ws.Cells.Item(292, 14).Value = "`$/unidad"
$ws.Cells.Item(292, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(292, 16).Value = 900
$ws.Cells.Item(292, 17).Value = 1
$ws.Cells.Item(292, 18).Value = "Hortaliza"
